# CUS15: actualización de servicios, scripts y archivos de cotización
#
# Adds a "Cotizacion N°" label + number to the quotation header (row 1,
# columns C/D) of the COTIZACION sheet, and updates the active selection
# to reflect the newly-entered cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COTIZACION")

# New label in C1, styled like the other section headers (e.g. A3
# "Datos del Cliente" / A9 "Datos del Emisor (Proveedor)") so it reuses
# the existing bold, dark-blue "label" look instead of inventing a new one.
$ws.Range("C1").Value = "Cotizacion N°"
$ws.Range("A3").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null

# New quotation number in D1 (plain number, default formatting).
$ws.Range("D1").Value = 102

# Reflect the edit in the sheet's active selection.
$ws.Range("C1:D1").Select() | Out-Null
